# Update cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.720.84'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.600.67'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '211.57'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("D10").Value = '19.66'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").Value = '1.825.73'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.585.54'
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '65.06'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '0.0₃0738'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '210.11'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '7.15'
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  -2.84%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '143.59'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("D28").Value = '15.35'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").Value = '2.97'
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").Value = '1.289.60'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -2.56%  '
$ws.Range("E37").Value = '  +10.71%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").Value = '0.831'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '5.39'
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").Value = '0.782'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").Value = '62.91'
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '1.737.25'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '90.52'
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '7.42'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("E51").Value = '  +0.96%  '
